$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").ClearContents()
$ws.Range("B2").Value = "보통이야"
$ws.Range("B3").Value = "바보야!"
$ws.Range("B4").Value = "보통이야"
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
